$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.77831899733424
$ws.Range("C2").Value = 10.27232470012285
$ws.Range("D2").Value = 14.41935727154581
$ws.Range("E2").Value = 15.4530111398583
$ws.Range("G2").Value = 3.666323054928875
$ws.Range("I2").Value = 23.26657578650508
$ws.Range("J2").Value = 9.104455859180565
$ws.Range("M2").Value = 18.70483167175254
$ws.Range("O2").Value = 26.66272790411299
$ws.Range("B3").Value = 14.24227262584655
$ws.Range("C3").Value = 9.8186725096886
$ws.Range("D3").Value = 14.41194776216358
$ws.Range("E3").Value = 15.48028740873552
$ws.Range("G3").Value = 3.668886759654261
$ws.Range("I3").Value = 23.40617755527639
$ws.Range("J3").Value = 9.128436658274119
$ws.Range("M3").Value = 18.53885790612103
$ws.Range("O3").Value = 26.75491709937675
$ws.Range("B4").Value = 13.90408216068654
$ws.Range("C4").Value = 9.52999519333731
$ws.Range("D4").Value = 14.41046869939988
$ws.Range("E4").Value = 15.49998531816587
$ws.Range("G4").Value = 3.670543684031637
$ws.Range("I4").Value = 23.49772831711715
$ws.Range("J4").Value = 9.14405141190778
$ws.Range("M4").Value = 18.4391115793822
$ws.Range("O4").Value = 26.81892833013275
$ws.Range("B5").Value = 13.76419235537879
$ws.Range("C5").Value = 9.4099659629385
$ws.Range("D5").Value = 14.4106391262835
$ws.Range("E5").Value = 15.50875319276621
$ws.Range("G5").Value = 3.671239784799373
$ws.Range("I5").Value = 23.53650048424958
$ws.Range("J5").Value = 9.150638921102138
$ws.Range("M5").Value = 18.39904415169284
$ws.Range("O5").Value = 26.84686786602606
$ws.Range("B6").Value = 13.74084479392143
$ws.Range("C6").Value = 9.389895820726638
$ws.Range("D6").Value = 14.41071414471945
$ws.Range("E6").Value = 15.51025380237388
$ws.Range("G6").Value = 3.671356635604667
$ws.Range("I6").Value = 23.54302694768378
$ws.Range("J6").Value = 9.151746336996499
$ws.Range("M6").Value = 18.39242705751563
$ws.Range("O6").Value = 26.85161896760209
$ws.Range("B7").Value = 13.90220367432801
$ws.Range("C7").Value = 9.528385894018539
$ws.Range("D7").Value = 14.41046786627479
$ws.Range("E7").Value = 15.50010056671898
$ws.Range("G7").Value = 3.670552987212295
$ws.Range("I7").Value = 23.49824528600789
$ws.Range("J7").Value = 9.144139344285032
$ws.Range("M7").Value = 18.43856882042701
$ws.Range("O7").Value = 26.8192976348291
$ws.Range("B8").Value = 14.59548507676061
$ws.Range("C8").Value = 10.11809621313642
$ws.Range("D8").Value = 14.41616599109942
$ws.Range("E8").Value = 15.46180311104399
$ws.Range("G8").Value = 3.667189877577528
$ws.Range("I8").Value = 23.31349783878481
$ws.Range("J8").Value = 9.11253990365525
$ws.Range("M8").Value = 18.64717721826413
$ws.Range("O8").Value = 26.69297273842657
$ws.Range("B9").Value = 15.87526870517945
$ws.Range("C9").Value = 11.18825888595637
$ws.Range("D9").Value = 14.45163653642799
$ws.Range("E9").Value = 15.41015160644133
$ws.Range("G9").Value = 3.661248606819189
$ws.Range("I9").Value = 22.99762742684996
$ws.Range("J9").Value = 9.057617909073743
$ws.Range("M9").Value = 19.07177730983194
$ws.Range("O9").Value = 26.50436525605669
$ws.Range("B10").Value = 16.75771414163906
$ws.Range("C10").Value = 11.91538398810733
$ws.Range("D10").Value = 14.49239339116539
$ws.Range("E10").Value = 15.38654860440219
$ws.Range("G10").Value = 3.657277636536972
$ws.Range("I10").Value = 22.79402786651801
$ws.Range("J10").Value = 9.021532157370729
$ws.Range("M10").Value = 19.39095430941829
$ws.Range("O10").Value = 26.40228782964145
$ws.Range("B11").Value = 17.14502391224839
$ws.Range("C11").Value = 12.23227686943509
$ws.Range("D11").Value = 14.51409150546555
$ws.Range("E11").Value = 15.37893347526451
$ws.Range("G11").Value = 3.655555756218159
$ws.Range("I11").Value = 22.70762682053729
$ws.Range("J11").Value = 9.006035917442155
$ws.Range("M11").Value = 19.53724032573222
$ws.Range("O11").Value = 26.36386397536291
$ws.Range("B12").Value = 17.28954732781705
$ws.Range("C12").Value = 12.35020986610091
$ws.Range("D12").Value = 14.52275831597074
$ws.Range("E12").Value = 15.37649906261644
$ws.Range("G12").Value = 3.654915808204836
$ws.Range("I12").Value = 22.67580670935577
$ws.Range("J12").Value = 9.000299645543729
$ws.Range("M12").Value = 19.59275102486056
$ws.Range("O12").Value = 26.35047265698749
$ws.Range("B13").Value = 17.25851853828891
$ws.Range("C13").Value = 12.32490385277361
$ws.Range("D13").Value = 14.52087180681572
$ws.Range("E13").Value = 15.37700337208915
$ws.Range("G13").Value = 3.655053095756136
$ws.Range("I13").Value = 22.68261973112245
$ws.Range("J13").Value = 9.001529197732431
$ws.Range("M13").Value = 19.58079133044552
$ws.Range("O13").Value = 26.35330506693696
$ws.Range("B14").Value = 17.15695746761226
$ws.Range("C14").Value = 12.24202109171323
$ws.Range("D14").Value = 14.51479553020621
$ws.Range("E14").Value = 15.37872418933589
$ws.Range("G14").Value = 3.655502865365904
$ws.Range("I14").Value = 22.70499093818576
$ws.Range("J14").Value = 9.005561351247875
$ws.Range("M14").Value = 19.54180514934787
$ws.Range("O14").Value = 26.36273899835894
$ws.Range("B15").Value = 17.09446629963943
$ws.Range("C15").Value = 12.19098183816156
$ws.Range("D15").Value = 14.51113214264077
$ws.Range("E15").Value = 15.37983675479646
$ws.Range("G15").Value = 3.655779934876685
$ws.Range("I15").Value = 22.71881102834969
$ws.Range("J15").Value = 9.008048317008885
$ws.Range("M15").Value = 19.51793878269632
$ws.Range("O15").Value = 26.36866868158793
$ws.Range("B16").Value = 16.73210850716052
$ws.Range("C16").Value = 11.89438860780893
$ws.Range("D16").Value = 14.491038589218
$ws.Range("E16").Value = 15.38710912602567
$ws.Range("G16").Value = 3.657391860922509
$ws.Range("I16").Value = 22.79979982991781
$ws.Range("J16").Value = 9.022563339821966
$ws.Range("M16").Value = 19.38141257292449
$ws.Range("O16").Value = 26.40496073364448
$ws.Range("B17").Value = 16.5061130247714
$ws.Range("C17").Value = 11.70882857471678
$ws.Range("D17").Value = 14.47951799818365
$ws.Range("E17").Value = 15.39237040125446
$ws.Range("G17").Value = 3.65840233000163
$ws.Range("I17").Value = 22.85107909101015
$ws.Range("J17").Value = 9.031703028674976
$ws.Range("M17").Value = 19.29790797573687
$ws.Range("O17").Value = 26.42928175951083
$ws.Range("B18").Value = 16.3748031584834
$ws.Range("C18").Value = 11.600796370886
$ws.Range("D18").Value = 14.47318913073673
$ws.Range("E18").Value = 15.39569039589247
$ws.Range("G18").Value = 3.658991485695879
$ws.Range("I18").Value = 22.8811585281691
$ws.Range("J18").Value = 9.037046491982242
$ws.Range("M18").Value = 19.24998349023076
$ws.Range("O18").Value = 26.44402428258672
$ws.Range("B19").Value = 16.33012028397006
$ws.Range("C19").Value = 11.56399710912753
$ws.Range("D19").Value = 14.47109748551905
$ws.Range("E19").Value = 15.39686494470751
$ws.Range("G19").Value = 3.659192332927241
$ws.Range("I19").Value = 22.89144325609199
$ws.Range("J19").Value = 9.038870574725388
$ws.Range("M19").Value = 19.2337764094167
$ws.Range("O19").Value = 26.44914508120588
$ws.Range("B20").Value = 16.5303084387481
$ws.Range("C20").Value = 11.72871716599377
$ws.Range("D20").Value = 14.48071362393543
$ws.Range("E20").Value = 15.39177991434877
$ws.Range("G20").Value = 3.658293940423275
$ws.Range("I20").Value = 22.84555975588045
$ws.Range("J20").Value = 9.03072113776453
$ws.Range("M20").Value = 19.30678660113349
$ws.Range("O20").Value = 26.42661468699804
$ws.Range("B21").Value = 17.18684738092483
$ws.Range("C21").Value = 12.26642239084357
$ws.Range("D21").Value = 14.5165680946296
$ws.Range("E21").Value = 15.37820654874032
$ws.Range("G21").Value = 3.655370429508118
$ws.Range("I21").Value = 22.69839556564621
$ws.Range("J21").Value = 9.004373435543602
$ws.Range("M21").Value = 19.55325353465043
$ws.Range("O21").Value = 26.35993651744791
$ws.Range("B22").Value = 17.60340404121911
$ws.Range("C22").Value = 12.60576776407457
$ws.Range("D22").Value = 14.54262307633539
$ws.Range("E22").Value = 15.37195429422778
$ws.Range("G22").Value = 3.653530189995449
$ws.Range("I22").Value = 22.60745198588623
$ws.Range("J22").Value = 8.987921851134125
$ws.Range("M22").Value = 19.71498726165213
$ws.Range("O22").Value = 26.32311666705266
$ws.Range("B23").Value = 17.38225969365984
$ws.Range("C23").Value = 12.42577842744806
$ws.Range("D23").Value = 14.52847852697955
$ws.Range("E23").Value = 15.37505156346672
$ws.Range("G23").Value = 3.654505936265681
$ws.Range("I23").Value = 22.65550975633219
$ws.Range("J23").Value = 8.996632209451009
$ws.Range("M23").Value = 19.62862079500736
$ws.Range("O23").Value = 26.34214758024035
$ws.Range("B24").Value = 16.51937398821426
$ws.Range("C24").Value = 11.7197297372847
$ws.Range("D24").Value = 14.48017216399428
$ws.Range("E24").Value = 15.39204595392331
$ws.Range("G24").Value = 3.658342917730048
$ws.Range("I24").Value = 22.84805318381346
$ws.Range("J24").Value = 9.031164773581681
$ws.Range("M24").Value = 19.30277231156713
$ws.Range("O24").Value = 26.42781810343616
$ws.Range("B25").Value = 15.53861168123829
$ws.Range("C25").Value = 10.90872704941715
$ws.Range("D25").Value = 14.4394496557327
$ws.Range("E25").Value = 15.42160848061653
$ws.Range("G25").Value = 3.662786351216485
$ws.Range("I25").Value = 23.07809199045041
$ws.Range("J25").Value = 9.071724635518901
$ws.Range("M25").Value = 18.95547929350261
$ws.Range("O25").Value = 26.54901155172463
